$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format for price cells whose new values would otherwise be
# auto-detected as numbers by Excel, then restore the default "Normal" style
# so only the value changes (matches original cell styling).
$textCells = @("D6", "D7", "D8", "D9", "D11", "D12", "D13", "D14", "D15", "D17", "D18", "D19", "D21", "D22", "D23", "D24", "D25", "D26", "D27", "D28", "D29", "D30", "D31", "D32", "D33", "D34", "D35", "D36", "D37", "D38", "D39", "D40", "D41", "D42", "D43", "D44", "D45", "D46", "D47", "D48", "D49", "D50", "D51")
foreach ($ref in $textCells) {
    $ws.Range($ref).NumberFormat = "@"
}

$ws.Range("D2").Value = "30.330.28"
$ws.Range("E2").Value = "  -0.15%  "
$ws.Range("D3").Value = "1.861.00"
$ws.Range("E3").Value = "  -0.72%  "
$ws.Range("E4").Value = "  +0.27%  "
$ws.Range("E5").Value = "  -2.44%  "
$ws.Range("D6").Value = "1.002"
$ws.Range("E6").Value = "  +0.22%  "
$ws.Range("D7").Value = "0.4751"
$ws.Range("E7").Value = "  -0.79%  "
$ws.Range("D8").Value = "0.2751"
$ws.Range("E8").Value = "  -2.79%  "
$ws.Range("D9").Value = "0.06445"
$ws.Range("E9").Value = "  -1.17%  "
$ws.Range("D10").Value = "1.861.09"
$ws.Range("E10").Value = "  -0.71%  "
$ws.Range("D11").Value = "0.07430"
$ws.Range("E11").Value = "  -0.35%  "
$ws.Range("D12").Value = "16.03"
$ws.Range("E12").Value = "  -3.41%  "
$ws.Range("D13").Value = "5.002"
$ws.Range("E13").Value = "  -1.87%  "
$ws.Range("D14").Value = "85.65"
$ws.Range("E14").Value = "  -2.94%  "
$ws.Range("D15").Value = "0.6306"
$ws.Range("E15").Value = "  -4.08%  "
$ws.Range("D16").Value = "30.301.98"
$ws.Range("E16").Value = "  -0.16%  "
$ws.Range("D17").Value = "1.002"
$ws.Range("E17").Value = "  +0.24%  "
$ws.Range("D18").Value = "12.82"
$ws.Range("E18").Value = "  -3.80%  "
$ws.Range("B19").Value = "ShibaInu"
$ws.Range("C19").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D19").Value = "0.000007330"
$ws.Range("E19").Value = "  -3.57%  "
$ws.Range("B20").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C20").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D20").Value = "2.114.02"
$ws.Range("E20").Value = "  -0.02%  "
$ws.Range("B21").Value = "BitcoinCash"
$ws.Range("C21").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D21").Value = "223.48"
$ws.Range("E21").Value = "  +2.01%  "
$ws.Range("B22").Value = "BinanceUSD"
$ws.Range("C22").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D22").Value = "1.003"
$ws.Range("E22").Value = "  +0.28%  "
$ws.Range("B23").Value = "Uniswap"
$ws.Range("C23").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D23").Value = "5.099"
$ws.Range("E23").Value = "  -4.08%  "
$ws.Range("B24").Value = "Chainlink"
$ws.Range("C24").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D24").Value = "6.048"
$ws.Range("E24").Value = "  -2.75%  "
$ws.Range("D25").Value = "9.220"
$ws.Range("E25").Value = "  -1.39%  "
$ws.Range("B26").Value = "Monero"
$ws.Range("C26").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D26").Value = "166.42"
$ws.Range("E26").Value = "  -0.71%  "
$ws.Range("B27").Value = "EthereumClassic"
$ws.Range("C27").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D27").Value = "17.75"
$ws.Range("E27").Value = "  -3.89%  "
$ws.Range("B28").Value = "LidoDAOToken"
$ws.Range("C28").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D28").Value = "1.872"
$ws.Range("E28").Value = "  -5.32%  "
$ws.Range("B29").Value = "Stellar"
$ws.Range("C29").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D29").Value = "0.1032"
$ws.Range("E29").Value = "  +9.41%  "
$ws.Range("B30").Value = "Toncoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D30").Value = "1.380"
$ws.Range("E30").Value = "  -5.50%  "
$ws.Range("B31").Value = "InternetComputer(DFINITY)"
$ws.Range("C31").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D31").Value = "4.209"
$ws.Range("E31").Value = "  -2.49%  "
$ws.Range("B32").Value = "Filecoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D32").Value = "3.883"
$ws.Range("E32").Value = "  -3.88%  "
$ws.Range("B33").Value = "Hedera"
$ws.Range("C33").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D33").Value = "0.04904"
$ws.Range("E33").Value = "  -3.19%  "
$ws.Range("B34").Value = "ARBITRUM"
$ws.Range("C34").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D34").Value = "1.153"
$ws.Range("E34").Value = "  -4.43%  "
$ws.Range("B35").Value = "ImmutableX"
$ws.Range("C35").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D35").Value = "0.7248"
$ws.Range("E35").Value = "  -3.50%  "
$ws.Range("B36").Value = "Frax"
$ws.Range("C36").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D36").Value = "1.002"
$ws.Range("E36").Value = "  +0.35%  "
$ws.Range("B37").Value = "HuobiToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D37").Value = "2.714"
$ws.Range("E37").Value = "  +0.09%  "
$ws.Range("B38").Value = "VeChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D38").Value = "0.01897"
$ws.Range("E38").Value = "  +3.96%  "
$ws.Range("B39").Value = "MXToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D39").Value = "2.627"
$ws.Range("E39").Value = "  +0.47%  "
$ws.Range("B40").Value = "TrustWalletToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D40").Value = "0.9037"
$ws.Range("E40").Value = "  -0.15%  "
$ws.Range("B41").Value = "RenderToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D41").Value = "1.981"
$ws.Range("E41").Value = "  -4.11%  "
$ws.Range("B42").Value = "Quant"
$ws.Range("C42").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D42").Value = "105.42"
$ws.Range("E42").Value = "  -1.40%  "
$ws.Range("B43").Value = "PaxDollar"
$ws.Range("C43").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D43").Value = "0.9960"
$ws.Range("E43").Value = "  -0.78%  "
$ws.Range("B44").Value = "TheSandbox"
$ws.Range("C44").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D44").Value = "0.4097"
$ws.Range("E44").Value = "  -4.17%  "
$ws.Range("B45").Value = "FraxShare"
$ws.Range("C45").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D45").Value = "5.552"
$ws.Range("E45").Value = "  -5.85%  "
$ws.Range("B46").Value = "Aptos"
$ws.Range("C46").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D46").Value = "7.035"
$ws.Range("E46").Value = "  -4.67%  "
$ws.Range("B47").Value = "Aave"
$ws.Range("C47").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D47").Value = "61.12"
$ws.Range("E47").Value = "  -5.15%  "
$ws.Range("B48").Value = "Algorand"
$ws.Range("C48").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D48").Value = "0.1206"
$ws.Range("E48").Value = "  -6.00%  "
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").Value = "8.789"
$ws.Range("E49").Value = "  -1.98%  "
$ws.Range("D50").Value = "1.398"
$ws.Range("E50").Value = "  -5.15%  "
$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D51").Value = "0.05609"
$ws.Range("E51").Value = "  -0.36%  "

foreach ($ref in $textCells) {
    $ws.Range($ref).Style = "Normal"
}